$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C values for existing rows 2-531 (new area projection values) ---
$ws.Cells.Item(2, 3).Value = 62920.84966366
$ws.Cells.Item(3, 3).Value = 62920.84966366
$ws.Cells.Item(4, 3).Value = 63214.57991599
$ws.Cells.Item(5, 3).Value = 63214.57991599
$ws.Cells.Item(6, 3).Value = 64689.18705694001
$ws.Cells.Item(7, 3).Value = 64689.18705694001
$ws.Cells.Item(8, 3).Value = 64806.21963948999
$ws.Cells.Item(9, 3).Value = 63867.56917078
$ws.Cells.Item(10, 3).Value = 64016.97279823
$ws.Cells.Item(11, 3).Value = 63966.75848065
$ws.Cells.Item(12, 3).Value = 63998.73993046
$ws.Cells.Item(13, 3).Value = 64186.56592389
$ws.Cells.Item(14, 3).Value = 64026.05337905
$ws.Cells.Item(15, 3).Value = 64026.05337905
$ws.Cells.Item(16, 3).Value = 63988.57673784999
$ws.Cells.Item(17, 3).Value = 63988.57673784999
$ws.Cells.Item(18, 3).Value = 64100.44749233
$ws.Cells.Item(19, 3).Value = 65375.29071184
$ws.Cells.Item(20, 3).Value = 65375.29071184
$ws.Cells.Item(21, 3).Value = 65375.29071184
$ws.Cells.Item(22, 3).Value = 66319.31041706
$ws.Cells.Item(23, 3).Value = 66319.31041706
$ws.Cells.Item(24, 3).Value = 66322.59855182
$ws.Cells.Item(25, 3).Value = 66119.08831507
$ws.Cells.Item(26, 3).Value = 66115.4281053
$ws.Cells.Item(27, 3).Value = 66115.4281053
$ws.Cells.Item(28, 3).Value = 66115.4281053
$ws.Cells.Item(29, 3).Value = 66149.19168074
$ws.Cells.Item(30, 3).Value = 66149.19168074
$ws.Cells.Item(31, 3).Value = 66149.19168074
$ws.Cells.Item(32, 3).Value = 66371.47592934
$ws.Cells.Item(33, 3).Value = 66360.6918424
$ws.Cells.Item(34, 3).Value = 68679.31219378
$ws.Cells.Item(35, 3).Value = 68609.33207737
$ws.Cells.Item(36, 3).Value = 71107.76774789
$ws.Cells.Item(37, 3).Value = 71349.56804622
$ws.Cells.Item(38, 3).Value = 71353.5277475
$ws.Cells.Item(39, 3).Value = 71353.5277475
$ws.Cells.Item(40, 3).Value = 71453.10991888
$ws.Cells.Item(41, 3).Value = 71453.15958974001
$ws.Cells.Item(42, 3).Value = 71457.01578338
$ws.Cells.Item(43, 3).Value = 71434.43294831
$ws.Cells.Item(44, 3).Value = 71433.73003696001
$ws.Cells.Item(45, 3).Value = 72091.04945063
$ws.Cells.Item(46, 3).Value = 72098.48507025
$ws.Cells.Item(47, 3).Value = 72098.48507025
$ws.Cells.Item(48, 3).Value = 72113.36945457001
$ws.Cells.Item(49, 3).Value = 72135.8306133
$ws.Cells.Item(50, 3).Value = 72147.95434294999
$ws.Cells.Item(51, 3).Value = 72153.32898214999
$ws.Cells.Item(52, 3).Value = 74173.75318869999
$ws.Cells.Item(53, 3).Value = 74390.90374638001
$ws.Cells.Item(54, 3).Value = 74446.39482011
$ws.Cells.Item(55, 3).Value = 76337.13687613999
$ws.Cells.Item(56, 3).Value = 76244.34752007
$ws.Cells.Item(57, 3).Value = 76263.43689525999
$ws.Cells.Item(58, 3).Value = 76314.72824959
$ws.Cells.Item(59, 3).Value = 76314.72824959
$ws.Cells.Item(60, 3).Value = 76520.42126783
$ws.Cells.Item(61, 3).Value = 76845.15176653999
$ws.Cells.Item(62, 3).Value = 76920.69111711001
$ws.Cells.Item(63, 3).Value = 77118.12046655999
$ws.Cells.Item(64, 3).Value = 77118.12928582
$ws.Cells.Item(65, 3).Value = 77161.94100288
$ws.Cells.Item(66, 3).Value = 77825.3239321
$ws.Cells.Item(67, 3).Value = 77893.35480553999
$ws.Cells.Item(68, 3).Value = 77885.4008187
$ws.Cells.Item(69, 3).Value = 77873.11221466999
$ws.Cells.Item(70, 3).Value = 77846.4635974
$ws.Cells.Item(71, 3).Value = 77843.08569771
$ws.Cells.Item(72, 3).Value = 77843.08569771
$ws.Cells.Item(73, 3).Value = 77843.09410647
$ws.Cells.Item(74, 3).Value = 77897.24196855001
$ws.Cells.Item(75, 3).Value = 77897.24196855001
$ws.Cells.Item(76, 3).Value = 77921.6790111
$ws.Cells.Item(77, 3).Value = 77935.53979914999
$ws.Cells.Item(78, 3).Value = 77925.28120869001
$ws.Cells.Item(79, 3).Value = 77925.28120869001
$ws.Cells.Item(80, 3).Value = 77948.6180761
$ws.Cells.Item(81, 3).Value = 77914.18755932001
$ws.Cells.Item(82, 3).Value = 77914.18755932001
$ws.Cells.Item(83, 3).Value = 78416.32795834
$ws.Cells.Item(84, 3).Value = 78416.32795834
$ws.Cells.Item(85, 3).Value = 78508.30801137
$ws.Cells.Item(86, 3).Value = 78501.59139307
$ws.Cells.Item(87, 3).Value = 78545.54305764
$ws.Cells.Item(88, 3).Value = 78544.45295228
$ws.Cells.Item(89, 3).Value = 78503.09296289999
$ws.Cells.Item(90, 3).Value = 78451.49910236
$ws.Cells.Item(91, 3).Value = 78470.03204619999
$ws.Cells.Item(92, 3).Value = 78470.03204619999
$ws.Cells.Item(93, 3).Value = 78283.1714218
$ws.Cells.Item(94, 3).Value = 78283.1714218
$ws.Cells.Item(95, 3).Value = 78328.54758742999
$ws.Cells.Item(96, 3).Value = 78442.25517402
$ws.Cells.Item(97, 3).Value = 78205.99407938001
$ws.Cells.Item(98, 3).Value = 78205.99407938001
$ws.Cells.Item(99, 3).Value = 78341.99207405999
$ws.Cells.Item(100, 3).Value = 78341.99207405999
$ws.Cells.Item(101, 3).Value = 78341.99207405999
$ws.Cells.Item(102, 3).Value = 78336.57724717999
$ws.Cells.Item(103, 3).Value = 78319.55157428001
$ws.Cells.Item(104, 3).Value = 78537.95088785001
$ws.Cells.Item(105, 3).Value = 121961.28043705
$ws.Cells.Item(106, 3).Value = 121962.53472412
$ws.Cells.Item(107, 3).Value = 121910.19427115
$ws.Cells.Item(108, 3).Value = 121924.94178551
$ws.Cells.Item(109, 3).Value = 121960.87149603
$ws.Cells.Item(110, 3).Value = 121957.69733438
$ws.Cells.Item(111, 3).Value = 121944.38149854
$ws.Cells.Item(112, 3).Value = 122169.27459826
$ws.Cells.Item(113, 3).Value = 121966.41465257
$ws.Cells.Item(114, 3).Value = 121838.60936998
$ws.Cells.Item(115, 3).Value = 121914.95144998
$ws.Cells.Item(116, 3).Value = 121928.39066953
$ws.Cells.Item(117, 3).Value = 121608.98426339
$ws.Cells.Item(118, 3).Value = 121687.61347036
$ws.Cells.Item(119, 3).Value = 122140.77398903
$ws.Cells.Item(120, 3).Value = 122213.2931441
$ws.Cells.Item(121, 3).Value = 122218.93458745
$ws.Cells.Item(122, 3).Value = 122266.14409314
$ws.Cells.Item(123, 3).Value = 122350.14599157
$ws.Cells.Item(124, 3).Value = 122319.50659697
$ws.Cells.Item(125, 3).Value = 122335.06548313
$ws.Cells.Item(126, 3).Value = 122360.94520154
$ws.Cells.Item(127, 3).Value = 122392.72664975
$ws.Cells.Item(128, 3).Value = 122537.94315933
$ws.Cells.Item(129, 3).Value = 122537.94315933
$ws.Cells.Item(130, 3).Value = 122589.51838731
$ws.Cells.Item(131, 3).Value = 122580.45674935
$ws.Cells.Item(132, 3).Value = 122593.29088971
$ws.Cells.Item(133, 3).Value = 122586.9764597
$ws.Cells.Item(134, 3).Value = 122653.70295204
$ws.Cells.Item(135, 3).Value = 122682.11896024
$ws.Cells.Item(136, 3).Value = 122718.91273881
$ws.Cells.Item(137, 3).Value = 122640.75153496
$ws.Cells.Item(138, 3).Value = 122693.01505771
$ws.Cells.Item(139, 3).Value = 122700.16051611
$ws.Cells.Item(140, 3).Value = 122675.52304192
$ws.Cells.Item(141, 3).Value = 122664.92340096
$ws.Cells.Item(142, 3).Value = 122637.77320595
$ws.Cells.Item(143, 3).Value = 122614.17041656
$ws.Cells.Item(144, 3).Value = 122572.69729503
$ws.Cells.Item(145, 3).Value = 122592.79101958
$ws.Cells.Item(146, 3).Value = 122685.94186011
$ws.Cells.Item(147, 3).Value = 122710.03042586
$ws.Cells.Item(148, 3).Value = 122710.03042586
$ws.Cells.Item(149, 3).Value = 122710.03042586
$ws.Cells.Item(150, 3).Value = 122712.77874165
$ws.Cells.Item(151, 3).Value = 122731.4900746
$ws.Cells.Item(152, 3).Value = 122731.4900746
$ws.Cells.Item(153, 3).Value = 122748.96138702
$ws.Cells.Item(154, 3).Value = 122762.63925898
$ws.Cells.Item(155, 3).Value = 122792.10227425
$ws.Cells.Item(156, 3).Value = 122724.40224028
$ws.Cells.Item(157, 3).Value = 122925.44402368
$ws.Cells.Item(158, 3).Value = 122919.4137614
$ws.Cells.Item(159, 3).Value = 122952.60867228
$ws.Cells.Item(160, 3).Value = 122952.60867228
$ws.Cells.Item(161, 3).Value = 123016.02977168
$ws.Cells.Item(162, 3).Value = 123036.27960343
$ws.Cells.Item(163, 3).Value = 122957.22590582
$ws.Cells.Item(164, 3).Value = 123008.80724765
$ws.Cells.Item(165, 3).Value = 123036.17246982
$ws.Cells.Item(166, 3).Value = 123036.17246982
$ws.Cells.Item(167, 3).Value = 123044.70973996
$ws.Cells.Item(168, 3).Value = 123129.98212332
$ws.Cells.Item(169, 3).Value = 123129.98212332
$ws.Cells.Item(170, 3).Value = 123124.62277556
$ws.Cells.Item(171, 3).Value = 123164.26319526
$ws.Cells.Item(172, 3).Value = 123172.85780978
$ws.Cells.Item(173, 3).Value = 123172.85780978
$ws.Cells.Item(174, 3).Value = 123172.85780978
$ws.Cells.Item(175, 3).Value = 123352.16094778
$ws.Cells.Item(176, 3).Value = 123501.95094421
$ws.Cells.Item(177, 3).Value = 123582.3810958
$ws.Cells.Item(178, 3).Value = 123598.58884438
$ws.Cells.Item(179, 3).Value = 123738.96260997
$ws.Cells.Item(180, 3).Value = 123755.14919889
$ws.Cells.Item(181, 3).Value = 123784.34384082
$ws.Cells.Item(182, 3).Value = 123784.34384082
$ws.Cells.Item(183, 3).Value = 123808.50279039
$ws.Cells.Item(184, 3).Value = 123814.34853426
$ws.Cells.Item(185, 3).Value = 123903.68146181
$ws.Cells.Item(186, 3).Value = 123938.2156508
$ws.Cells.Item(187, 3).Value = 123949.60760476
$ws.Cells.Item(188, 3).Value = 123950.94455655
$ws.Cells.Item(189, 3).Value = 123965.1650036
$ws.Cells.Item(190, 3).Value = 123966.62363682
$ws.Cells.Item(191, 3).Value = 123981.11704026
$ws.Cells.Item(192, 3).Value = 123879.86163471
$ws.Cells.Item(193, 3).Value = 123894.45074888
$ws.Cells.Item(194, 3).Value = 123939.01545252
$ws.Cells.Item(195, 3).Value = 124001.71922285
$ws.Cells.Item(196, 3).Value = 124016.9747788
$ws.Cells.Item(197, 3).Value = 124022.67117416
$ws.Cells.Item(198, 3).Value = 124016.85942214
$ws.Cells.Item(199, 3).Value = 124036.38405172
$ws.Cells.Item(200, 3).Value = 124037.03574398
$ws.Cells.Item(201, 3).Value = 124122.63648624
$ws.Cells.Item(202, 3).Value = 124157.77351202
$ws.Cells.Item(203, 3).Value = 124168.89775803
$ws.Cells.Item(204, 3).Value = 124194.45745338
$ws.Cells.Item(205, 3).Value = 124193.51472874
$ws.Cells.Item(206, 3).Value = 124252.08248889
$ws.Cells.Item(207, 3).Value = 124239.08045675
$ws.Cells.Item(208, 3).Value = 124309.64235124
$ws.Cells.Item(209, 3).Value = 124318.95207719
$ws.Cells.Item(210, 3).Value = 124326.09081902
$ws.Cells.Item(211, 3).Value = 124354.62187015
$ws.Cells.Item(212, 3).Value = 124344.80646756
$ws.Cells.Item(213, 3).Value = 124434.07555947
$ws.Cells.Item(214, 3).Value = 124434.07555947
$ws.Cells.Item(215, 3).Value = 124434.07555947
$ws.Cells.Item(216, 3).Value = 124562.99545036
$ws.Cells.Item(217, 3).Value = 124563.01289762
$ws.Cells.Item(218, 3).Value = 124563.01289762
$ws.Cells.Item(219, 3).Value = 124820.11622199
$ws.Cells.Item(220, 3).Value = 124849.89603684
$ws.Cells.Item(221, 3).Value = 124857.62131918
$ws.Cells.Item(222, 3).Value = 124858.09719745
$ws.Cells.Item(223, 3).Value = 124864.92346169
$ws.Cells.Item(224, 3).Value = 124881.77128389
$ws.Cells.Item(225, 3).Value = 124869.62395964
$ws.Cells.Item(226, 3).Value = 124931.0273543
$ws.Cells.Item(227, 3).Value = 124929.00479073
$ws.Cells.Item(228, 3).Value = 124920.7855369
$ws.Cells.Item(229, 3).Value = 124934.04758649
$ws.Cells.Item(230, 3).Value = 124931.76264068
$ws.Cells.Item(231, 3).Value = 124896.51730971
$ws.Cells.Item(232, 3).Value = 124902.39601054
$ws.Cells.Item(233, 3).Value = 124913.791162
$ws.Cells.Item(234, 3).Value = 124890.38969698
$ws.Cells.Item(235, 3).Value = 124783.89114681
$ws.Cells.Item(236, 3).Value = 124822.0850154
$ws.Cells.Item(237, 3).Value = 124822.0850154
$ws.Cells.Item(238, 3).Value = 124827.00960901
$ws.Cells.Item(239, 3).Value = 124813.55735009
$ws.Cells.Item(240, 3).Value = 124817.69706978
$ws.Cells.Item(241, 3).Value = 124783.34865494
$ws.Cells.Item(242, 3).Value = 124789.04553201
$ws.Cells.Item(243, 3).Value = 124832.59073303
$ws.Cells.Item(244, 3).Value = 124830.98823461
$ws.Cells.Item(245, 3).Value = 124830.98823461
$ws.Cells.Item(246, 3).Value = 124805.34252807
$ws.Cells.Item(247, 3).Value = 124824.58399574
$ws.Cells.Item(248, 3).Value = 124835.41067046
$ws.Cells.Item(249, 3).Value = 124856.79831218
$ws.Cells.Item(250, 3).Value = 124856.79831218
$ws.Cells.Item(251, 3).Value = 124856.79831218
$ws.Cells.Item(252, 3).Value = 124844.86161089
$ws.Cells.Item(253, 3).Value = 124840.94153193
$ws.Cells.Item(254, 3).Value = 124835.24467512
$ws.Cells.Item(255, 3).Value = 124824.26613292
$ws.Cells.Item(256, 3).Value = 124839.45124586
$ws.Cells.Item(257, 3).Value = 124859.90285494
$ws.Cells.Item(258, 3).Value = 124898.89923253
$ws.Cells.Item(259, 3).Value = 124850.8309061
$ws.Cells.Item(260, 3).Value = 124833.68277911
$ws.Cells.Item(261, 3).Value = 124817.81016734
$ws.Cells.Item(262, 3).Value = 124812.70856239
$ws.Cells.Item(263, 3).Value = 124831.44319969
$ws.Cells.Item(264, 3).Value = 124836.48422576
$ws.Cells.Item(265, 3).Value = 124866.45734432
$ws.Cells.Item(266, 3).Value = 124883.8890357
$ws.Cells.Item(267, 3).Value = 124888.42944423
$ws.Cells.Item(268, 3).Value = 124888.42944423
$ws.Cells.Item(269, 3).Value = 124917.05770224
$ws.Cells.Item(270, 3).Value = 124919.03798794
$ws.Cells.Item(271, 3).Value = 124927.22336394
$ws.Cells.Item(272, 3).Value = 124927.22336394
$ws.Cells.Item(273, 3).Value = 124931.53671034
$ws.Cells.Item(274, 3).Value = 125041.11968798
$ws.Cells.Item(275, 3).Value = 125055.77558577
$ws.Cells.Item(276, 3).Value = 125088.37680808
$ws.Cells.Item(277, 3).Value = 125088.37680808
$ws.Cells.Item(278, 3).Value = 125090.31324916
$ws.Cells.Item(279, 3).Value = 125090.31324916
$ws.Cells.Item(280, 3).Value = 125134.70954348
$ws.Cells.Item(281, 3).Value = 125141.19216682
$ws.Cells.Item(282, 3).Value = 125142.19173708
$ws.Cells.Item(283, 3).Value = 125140.50437506
$ws.Cells.Item(284, 3).Value = 125189.3823319
$ws.Cells.Item(285, 3).Value = 125198.31655446
$ws.Cells.Item(286, 3).Value = 125179.7181336
$ws.Cells.Item(287, 3).Value = 125179.7181336
$ws.Cells.Item(288, 3).Value = 125177.72603196
$ws.Cells.Item(289, 3).Value = 125179.25069799
$ws.Cells.Item(290, 3).Value = 125092.72131696
$ws.Cells.Item(291, 3).Value = 125080.29260493
$ws.Cells.Item(292, 3).Value = 124997.38030809
$ws.Cells.Item(293, 3).Value = 124844.70784338
$ws.Cells.Item(294, 3).Value = 124724.39461811
$ws.Cells.Item(295, 3).Value = 124690.03044382
$ws.Cells.Item(296, 3).Value = 124657.59000282
$ws.Cells.Item(297, 3).Value = 124418.20356032
$ws.Cells.Item(298, 3).Value = 124227.60355634
$ws.Cells.Item(299, 3).Value = 124196.17761418
$ws.Cells.Item(300, 3).Value = 124034.1450489
$ws.Cells.Item(301, 3).Value = 123691.6599318
$ws.Cells.Item(302, 3).Value = 123691.6599318
$ws.Cells.Item(303, 3).Value = 123691.6599318
$ws.Cells.Item(304, 3).Value = 123468.15076206
$ws.Cells.Item(305, 3).Value = 123232.8533199
$ws.Cells.Item(306, 3).Value = 122487.16957234
$ws.Cells.Item(307, 3).Value = 120814.33069119
$ws.Cells.Item(308, 3).Value = 120382.91425095
$ws.Cells.Item(309, 3).Value = 119693.61714513
$ws.Cells.Item(310, 3).Value = 119002.9437661
$ws.Cells.Item(311, 3).Value = 117235.68952503
$ws.Cells.Item(312, 3).Value = 116238.18093334
$ws.Cells.Item(313, 3).Value = 116148.11418568
$ws.Cells.Item(314, 3).Value = 115799.16298925
$ws.Cells.Item(315, 3).Value = 115766.19794197
$ws.Cells.Item(316, 3).Value = 115782.22732114
$ws.Cells.Item(317, 3).Value = 115671.82685922
$ws.Cells.Item(318, 3).Value = 115839.6022011
$ws.Cells.Item(319, 3).Value = 115823.21311319
$ws.Cells.Item(320, 3).Value = 115739.03615569
$ws.Cells.Item(321, 3).Value = 115667.63458987
$ws.Cells.Item(322, 3).Value = 116534.38493613
$ws.Cells.Item(323, 3).Value = 116534.38493613
$ws.Cells.Item(324, 3).Value = 116510.76247895
$ws.Cells.Item(325, 3).Value = 116468.73374515
$ws.Cells.Item(326, 3).Value = 116476.5396342
$ws.Cells.Item(327, 3).Value = 116431.85835027
$ws.Cells.Item(328, 3).Value = 116741.05660201
$ws.Cells.Item(329, 3).Value = 116667.18320714
$ws.Cells.Item(330, 3).Value = 116667.18320714
$ws.Cells.Item(331, 3).Value = 116598.97545711
$ws.Cells.Item(332, 3).Value = 116615.89100186
$ws.Cells.Item(333, 3).Value = 116533.63734644
$ws.Cells.Item(334, 3).Value = 116340.15991904
$ws.Cells.Item(335, 3).Value = 116337.99818145
$ws.Cells.Item(336, 3).Value = 116227.07810773
$ws.Cells.Item(337, 3).Value = 116133.63105732
$ws.Cells.Item(338, 3).Value = 116133.63105732
$ws.Cells.Item(339, 3).Value = 115831.73148276
$ws.Cells.Item(340, 3).Value = 115601.92505959
$ws.Cells.Item(341, 3).Value = 115598.59352755
$ws.Cells.Item(342, 3).Value = 115546.21657178
$ws.Cells.Item(343, 3).Value = 115506.57871061
$ws.Cells.Item(344, 3).Value = 115420.44881361
$ws.Cells.Item(345, 3).Value = 115235.81683435
$ws.Cells.Item(346, 3).Value = 115219.51881977
$ws.Cells.Item(347, 3).Value = 115215.77592972
$ws.Cells.Item(348, 3).Value = 115071.85258577
$ws.Cells.Item(349, 3).Value = 114875.88096072
$ws.Cells.Item(350, 3).Value = 114546.14652338
$ws.Cells.Item(351, 3).Value = 114546.14652338
$ws.Cells.Item(352, 3).Value = 114281.24434674
$ws.Cells.Item(353, 3).Value = 113831.48799207
$ws.Cells.Item(354, 3).Value = 113831.48799207
$ws.Cells.Item(355, 3).Value = 113742.34262015
$ws.Cells.Item(356, 3).Value = 113560.79259307
$ws.Cells.Item(357, 3).Value = 113379.36898293
$ws.Cells.Item(358, 3).Value = 113379.36898293
$ws.Cells.Item(359, 3).Value = 113040.51874646
$ws.Cells.Item(360, 3).Value = 113069.5805629
$ws.Cells.Item(361, 3).Value = 113069.5805629
$ws.Cells.Item(362, 3).Value = 113070.86179032
$ws.Cells.Item(363, 3).Value = 113130.54043847
$ws.Cells.Item(364, 3).Value = 113117.95325513
$ws.Cells.Item(365, 3).Value = 113148.50540621
$ws.Cells.Item(366, 3).Value = 113148.50540621
$ws.Cells.Item(367, 3).Value = 113027.68732558
$ws.Cells.Item(368, 3).Value = 113022.18652121
$ws.Cells.Item(369, 3).Value = 112966.86516115
$ws.Cells.Item(370, 3).Value = 113015.91748001
$ws.Cells.Item(371, 3).Value = 113011.25893998
$ws.Cells.Item(372, 3).Value = 113153.95717747
$ws.Cells.Item(373, 3).Value = 113137.26638155
$ws.Cells.Item(374, 3).Value = 113136.78168573
$ws.Cells.Item(375, 3).Value = 113143.15598394
$ws.Cells.Item(376, 3).Value = 113246.04543436
$ws.Cells.Item(377, 3).Value = 113304.44798242
$ws.Cells.Item(378, 3).Value = 113346.81620964
$ws.Cells.Item(379, 3).Value = 113329.49598246
$ws.Cells.Item(380, 3).Value = 113205.7372292
$ws.Cells.Item(381, 3).Value = 113216.82868722
$ws.Cells.Item(382, 3).Value = 113248.85630082
$ws.Cells.Item(383, 3).Value = 113220.25491014
$ws.Cells.Item(384, 3).Value = 113218.69182342
$ws.Cells.Item(385, 3).Value = 113221.50162936
$ws.Cells.Item(386, 3).Value = 113193.74381843
$ws.Cells.Item(387, 3).Value = 113212.98089333
$ws.Cells.Item(388, 3).Value = 113274.26038529
$ws.Cells.Item(389, 3).Value = 113274.26038529
$ws.Cells.Item(390, 3).Value = 113268.52434887
$ws.Cells.Item(391, 3).Value = 113268.52434887
$ws.Cells.Item(392, 3).Value = 113263.10121803
$ws.Cells.Item(393, 3).Value = 113259.98563409
$ws.Cells.Item(394, 3).Value = 113261.49662022
$ws.Cells.Item(395, 3).Value = 113283.04453747
$ws.Cells.Item(396, 3).Value = 113265.40012513
$ws.Cells.Item(397, 3).Value = 113114.10550717
$ws.Cells.Item(398, 3).Value = 113057.31587493
$ws.Cells.Item(399, 3).Value = 112914.15136645
$ws.Cells.Item(400, 3).Value = 112514.4442759
$ws.Cells.Item(401, 3).Value = 112458.52835665
$ws.Cells.Item(402, 3).Value = 111915.73673453
$ws.Cells.Item(403, 3).Value = 111381.25397377
$ws.Cells.Item(404, 3).Value = 111111.51666858
$ws.Cells.Item(405, 3).Value = 109592.41819814
$ws.Cells.Item(406, 3).Value = 109239.32247663
$ws.Cells.Item(407, 3).Value = 108698.33847377
$ws.Cells.Item(408, 3).Value = 108702.95001848
$ws.Cells.Item(409, 3).Value = 108697.78900613
$ws.Cells.Item(410, 3).Value = 108673.29704941
$ws.Cells.Item(411, 3).Value = 108673.29704941
$ws.Cells.Item(412, 3).Value = 108682.44353315
$ws.Cells.Item(413, 3).Value = 108682.87489085
$ws.Cells.Item(414, 3).Value = 108680.98093892
$ws.Cells.Item(415, 3).Value = 108676.33157682
$ws.Cells.Item(416, 3).Value = 108678.77746759
$ws.Cells.Item(417, 3).Value = 108686.86716162
$ws.Cells.Item(418, 3).Value = 108448.15878296
$ws.Cells.Item(419, 3).Value = 108453.44290499
$ws.Cells.Item(420, 3).Value = 108534.57836023
$ws.Cells.Item(421, 3).Value = 108534.57836023
$ws.Cells.Item(422, 3).Value = 108537.27909687
$ws.Cells.Item(423, 3).Value = 108546.18740486
$ws.Cells.Item(424, 3).Value = 108558.18503729
$ws.Cells.Item(425, 3).Value = 108560.55165108
$ws.Cells.Item(426, 3).Value = 108523.41650927
$ws.Cells.Item(427, 3).Value = 108531.79452511
$ws.Cells.Item(428, 3).Value = 108591.67695408
$ws.Cells.Item(429, 3).Value = 108592.57730295
$ws.Cells.Item(430, 3).Value = 108591.89301205
$ws.Cells.Item(431, 3).Value = 108617.38159267
$ws.Cells.Item(432, 3).Value = 108625.08741793
$ws.Cells.Item(433, 3).Value = 108626.80586908
$ws.Cells.Item(434, 3).Value = 108628.9026975
$ws.Cells.Item(435, 3).Value = 108630.89151317
$ws.Cells.Item(436, 3).Value = 108629.47034503
$ws.Cells.Item(437, 3).Value = 108617.28558268
$ws.Cells.Item(438, 3).Value = 108669.38740316
$ws.Cells.Item(439, 3).Value = 108680.19081614
$ws.Cells.Item(440, 3).Value = 108681.49564846
$ws.Cells.Item(441, 3).Value = 108682.59912121
$ws.Cells.Item(442, 3).Value = 108682.97895581
$ws.Cells.Item(443, 3).Value = 108683.6307655
$ws.Cells.Item(444, 3).Value = 108684.57459476
$ws.Cells.Item(445, 3).Value = 108684.57459476
$ws.Cells.Item(446, 3).Value = 108685.71145137
$ws.Cells.Item(447, 3).Value = 108689.92527425
$ws.Cells.Item(448, 3).Value = 108702.28238008
$ws.Cells.Item(449, 3).Value = 108702.52751655
$ws.Cells.Item(450, 3).Value = 108713.96073111
$ws.Cells.Item(451, 3).Value = 108722.30114313
$ws.Cells.Item(452, 3).Value = 108722.30114313
$ws.Cells.Item(453, 3).Value = 108723.82967512
$ws.Cells.Item(454, 3).Value = 108723.2920416
$ws.Cells.Item(455, 3).Value = 108723.2920416
$ws.Cells.Item(456, 3).Value = 108718.40691596
$ws.Cells.Item(457, 3).Value = 108722.04137561
$ws.Cells.Item(458, 3).Value = 108724.49470255
$ws.Cells.Item(459, 3).Value = 108723.94825406
$ws.Cells.Item(460, 3).Value = 108749.03488553
$ws.Cells.Item(461, 3).Value = 108749.14226977
$ws.Cells.Item(462, 3).Value = 108755.88937851
$ws.Cells.Item(463, 3).Value = 108770.57313768
$ws.Cells.Item(464, 3).Value = 108762.64541579
$ws.Cells.Item(465, 3).Value = 108772.49453963
$ws.Cells.Item(466, 3).Value = 108781.24495701
$ws.Cells.Item(467, 3).Value = 108790.08490544
$ws.Cells.Item(468, 3).Value = 108789.76048755
$ws.Cells.Item(469, 3).Value = 108784.91554714
$ws.Cells.Item(470, 3).Value = 108784.91554714
$ws.Cells.Item(471, 3).Value = 108784.91554714
$ws.Cells.Item(472, 3).Value = 108785.2915559
$ws.Cells.Item(473, 3).Value = 108833.27472732
$ws.Cells.Item(474, 3).Value = 108833.72079142
$ws.Cells.Item(475, 3).Value = 108836.78986723
$ws.Cells.Item(476, 3).Value = 108836.9175853
$ws.Cells.Item(477, 3).Value = 108881.40199022
$ws.Cells.Item(478, 3).Value = 108862.50139058
$ws.Cells.Item(479, 3).Value = 108862.50139058
$ws.Cells.Item(480, 3).Value = 108869.43422765
$ws.Cells.Item(481, 3).Value = 108870.96988734
$ws.Cells.Item(482, 3).Value = 108875.43642171
$ws.Cells.Item(483, 3).Value = 108883.07930522
$ws.Cells.Item(484, 3).Value = 108898.38025272
$ws.Cells.Item(485, 3).Value = 108901.0704302
$ws.Cells.Item(486, 3).Value = 108900.68172438
$ws.Cells.Item(487, 3).Value = 108900.68172438
$ws.Cells.Item(488, 3).Value = 108902.72156459
$ws.Cells.Item(489, 3).Value = 108905.88759933
$ws.Cells.Item(490, 3).Value = 108918.24424234
$ws.Cells.Item(491, 3).Value = 108918.24424234
$ws.Cells.Item(492, 3).Value = 108923.87408837
$ws.Cells.Item(493, 3).Value = 108930.74111592
$ws.Cells.Item(494, 3).Value = 108943.0469027
$ws.Cells.Item(495, 3).Value = 108943.07138778
$ws.Cells.Item(496, 3).Value = 108945.26065161
$ws.Cells.Item(497, 3).Value = 108953.15641341
$ws.Cells.Item(498, 3).Value = 108968.57761536
$ws.Cells.Item(499, 3).Value = 108975.38776857
$ws.Cells.Item(500, 3).Value = 108978.42741255
$ws.Cells.Item(501, 3).Value = 108978.95279825
$ws.Cells.Item(502, 3).Value = 108979.09543129
$ws.Cells.Item(503, 3).Value = 108983.32012715
$ws.Cells.Item(504, 3).Value = 109008.73801967
$ws.Cells.Item(505, 3).Value = 109022.02499029
$ws.Cells.Item(506, 3).Value = 109025.88959254
$ws.Cells.Item(507, 3).Value = 109024.59805526
$ws.Cells.Item(508, 3).Value = 109027.29709849
$ws.Cells.Item(509, 3).Value = 109027.87908413
$ws.Cells.Item(510, 3).Value = 109031.94039429
$ws.Cells.Item(511, 3).Value = 109035.79046511
$ws.Cells.Item(512, 3).Value = 109048.93534851
$ws.Cells.Item(513, 3).Value = 109067.82129496
$ws.Cells.Item(514, 3).Value = 109084.87642944
$ws.Cells.Item(515, 3).Value = 109097.08051312
$ws.Cells.Item(516, 3).Value = 109096.96209519
$ws.Cells.Item(517, 3).Value = 109101.48523202
$ws.Cells.Item(518, 3).Value = 109107.21215637
$ws.Cells.Item(519, 3).Value = 109110.67118091
$ws.Cells.Item(520, 3).Value = 109109.20219063
$ws.Cells.Item(521, 3).Value = 109109.51175791
$ws.Cells.Item(522, 3).Value = 109109.69944321
$ws.Cells.Item(523, 3).Value = 109112.58601586
$ws.Cells.Item(524, 3).Value = 109119.9922244
$ws.Cells.Item(525, 3).Value = 109127.24998387
$ws.Cells.Item(526, 3).Value = 109124.69313171
$ws.Cells.Item(527, 3).Value = 109130.13976517
$ws.Cells.Item(528, 3).Value = 109140.11003654
$ws.Cells.Item(529, 3).Value = 109142.96072876
$ws.Cells.Item(530, 3).Value = 109142.86240892
$ws.Cells.Item(531, 3).Value = 109145.10674065

# --- Append new rows 532-547 ---
# row 532
$ws.Range("A531:C531").Copy($ws.Range("A532:C532"))
$ws.Cells.Item(532, 1).Value = 44988.00577546296
$ws.Cells.Item(532, 2).NumberFormat = "@"
$ws.Cells.Item(532, 2).Value = "1677798499"
$ws.Cells.Item(532, 2).NumberFormat = "General"
$ws.Cells.Item(532, 3).Value = 109144.22581393
# row 533
$ws.Range("A531:C531").Copy($ws.Range("A533:C533"))
$ws.Cells.Item(533, 1).Value = 44988.01210648148
$ws.Cells.Item(533, 2).NumberFormat = "@"
$ws.Cells.Item(533, 2).Value = "1677799046"
$ws.Cells.Item(533, 2).NumberFormat = "General"
$ws.Cells.Item(533, 3).Value = 109144.22581393
# row 534
$ws.Range("A531:C531").Copy($ws.Range("A534:C534"))
$ws.Cells.Item(534, 1).Value = 44989.03908564815
$ws.Cells.Item(534, 2).NumberFormat = "@"
$ws.Cells.Item(534, 2).Value = "1677887777"
$ws.Cells.Item(534, 2).NumberFormat = "General"
$ws.Cells.Item(534, 3).Value = 109144.86160044
# row 535
$ws.Range("A531:C531").Copy($ws.Range("A535:C535"))
$ws.Cells.Item(535, 1).Value = 44990.02949074074
$ws.Cells.Item(535, 2).NumberFormat = "@"
$ws.Cells.Item(535, 2).Value = "1677973348"
$ws.Cells.Item(535, 2).NumberFormat = "General"
$ws.Cells.Item(535, 3).Value = 109145.94831197
# row 536
$ws.Range("A531:C531").Copy($ws.Range("A536:C536"))
$ws.Cells.Item(536, 1).Value = 44991.01837962963
$ws.Cells.Item(536, 2).NumberFormat = "@"
$ws.Cells.Item(536, 2).Value = "1678058788"
$ws.Cells.Item(536, 2).NumberFormat = "General"
$ws.Cells.Item(536, 3).Value = 109147.3061139
# row 537
$ws.Range("A531:C531").Copy($ws.Range("A537:C537"))
$ws.Cells.Item(537, 1).Value = 44992.1540625
$ws.Cells.Item(537, 2).NumberFormat = "@"
$ws.Cells.Item(537, 2).Value = "1678156911"
$ws.Cells.Item(537, 2).NumberFormat = "General"
$ws.Cells.Item(537, 3).Value = 109157.97074703
# row 538
$ws.Range("A531:C531").Copy($ws.Range("A538:C538"))
$ws.Cells.Item(538, 1).Value = 44993.03241898148
$ws.Cells.Item(538, 2).NumberFormat = "@"
$ws.Cells.Item(538, 2).Value = "1678232801"
$ws.Cells.Item(538, 2).NumberFormat = "General"
$ws.Cells.Item(538, 3).Value = 109166.9156174
# row 539
$ws.Range("A531:C531").Copy($ws.Range("A539:C539"))
$ws.Cells.Item(539, 1).Value = 44994.02148148148
$ws.Cells.Item(539, 2).NumberFormat = "@"
$ws.Cells.Item(539, 2).Value = "1678318256"
$ws.Cells.Item(539, 2).NumberFormat = "General"
$ws.Cells.Item(539, 3).Value = 109166.9156174
# row 540
$ws.Range("A531:C531").Copy($ws.Range("A540:C540"))
$ws.Cells.Item(540, 1).Value = 44995.39335648148
$ws.Cells.Item(540, 2).NumberFormat = "@"
$ws.Cells.Item(540, 2).Value = "1678436786"
$ws.Cells.Item(540, 2).NumberFormat = "General"
$ws.Cells.Item(540, 3).Value = 109172.68510872
# row 541
$ws.Range("A531:C531").Copy($ws.Range("A541:C541"))
$ws.Cells.Item(541, 1).Value = 44995.96583333334
$ws.Cells.Item(541, 2).NumberFormat = "@"
$ws.Cells.Item(541, 2).Value = "1678486248"
$ws.Cells.Item(541, 2).NumberFormat = "General"
$ws.Cells.Item(541, 3).Value = 109174.59736808
# row 542
$ws.Range("A531:C531").Copy($ws.Range("A542:C542"))
$ws.Cells.Item(542, 1).Value = 44997.02534722222
$ws.Cells.Item(542, 2).NumberFormat = "@"
$ws.Cells.Item(542, 2).Value = "1678577790"
$ws.Cells.Item(542, 2).NumberFormat = "General"
$ws.Cells.Item(542, 3).Value = 109179.43015703
# row 543
$ws.Range("A531:C531").Copy($ws.Range("A543:C543"))
$ws.Cells.Item(543, 1).Value = 44998.01040509259
$ws.Cells.Item(543, 2).NumberFormat = "@"
$ws.Cells.Item(543, 2).Value = "1678662899"
$ws.Cells.Item(543, 2).NumberFormat = "General"
$ws.Cells.Item(543, 3).Value = 109191.91109749
# row 544
$ws.Range("A531:C531").Copy($ws.Range("A544:C544"))
$ws.Cells.Item(544, 1).Value = 44999.03914351852
$ws.Cells.Item(544, 2).NumberFormat = "@"
$ws.Cells.Item(544, 2).Value = "1678751782"
$ws.Cells.Item(544, 2).NumberFormat = "General"
$ws.Cells.Item(544, 3).Value = 109191.54982029
# row 545
$ws.Range("A531:C531").Copy($ws.Range("A545:C545"))
$ws.Cells.Item(545, 1).Value = 44999.971875
$ws.Cells.Item(545, 2).NumberFormat = "@"
$ws.Cells.Item(545, 2).Value = "1678832370"
$ws.Cells.Item(545, 2).NumberFormat = "General"
$ws.Cells.Item(545, 3).Value = 109193.85578451
# row 546
$ws.Range("A531:C531").Copy($ws.Range("A546:C546"))
$ws.Cells.Item(546, 1).Value = 45001.04237268519
$ws.Cells.Item(546, 2).NumberFormat = "@"
$ws.Cells.Item(546, 2).Value = "1678924861"
$ws.Cells.Item(546, 2).NumberFormat = "General"
$ws.Cells.Item(546, 3).Value = 109194.38292894
# row 547
$ws.Range("A531:C531").Copy($ws.Range("A547:C547"))
$ws.Cells.Item(547, 1).Value = 45001.08371527777
$ws.Cells.Item(547, 2).NumberFormat = "@"
$ws.Cells.Item(547, 2).Value = "1678928433"
$ws.Cells.Item(547, 2).NumberFormat = "General"
$ws.Cells.Item(547, 3).Value = 109194.33406936
